$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data (row 15)
$ws.Range("A15").Value = "security/pgp/runme.sh"
$ws.Range("B15").Value = 0.04
$ws.Range("C15").Value = 0.01
$ws.Range("D15").Value = 0.02

# Match the selection state recorded in the saved file
$ws.Range("A32").Select()
